$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109864/1750183612496/A5MTrC4t3W-camera-video-segment-driver-1750183614996.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSBYDP2TS4%2F20250618%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250618T210004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKr%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIFpWzNZHTIuWx9Jfzil6EW78cqt8sJ%2FyIKXs7AqYf58IAiEA%2FplU9EnLA%2FBKbqOQgPr%2BpbyuqGuuuXTTshUMwZx9cVMq5gMIk%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDGPA%2FyuA9QEJpPr9DCq6Azw6oOyCKDKWrzL%2BHT%2FvK6UaJ2Lzu5Cf4R24oyE%2FlinXD0XsbTuWcmyYsn33d0ziGODpr1IxKv2fGo8pxPMoRS88PhZQOKtx3bXwg9nww1G8NH5neYXf16JZbDbTVVoU%2FAO5jWTR%2BR0%2FscWsGnNq%2FDPDnG%2FDp4ZD4o47mwcRExrnTCRFKs3Bjoi9ZRlko16yUzNitA8Pc8TPl3XgutBSpzQeAnADoQC44bB512exlSR8nRx6vJwvniQv3EAizZwtO8xT8oW%2BU3GkpXLOEa1Ql7u0HXNTl7WRSUm0ZqsLkGHKwxkwkPNy9fCheupmAAewRmeUMF6J4KnidkaH8WdLTp9urm1cFN4Kuutl2npiaOS%2FnhxV%2F3A4SN7XjRZKX89F3uytxQnl%2BjFW9%2BiK%2Ba6tBb34%2FM95IEXlATDB1rHKO5rnuhzqk784KOb7KXNUs76pLZasYyBFfU6s7wVXlO8VQdVI8I56j9JDuOFwqD9Vc%2BpRCp61IkkIp%2BQVIrS0l9VyiFge3JRBXfhU%2BwW94c2VOOSZkCGKotw%2Bp8uzJ9uaXeGYipFxMoLH3JHNgF1BKE42l2TWGb8gjLZBFuswtffLwgY6pQFdUfnEVvfK7XJ%2FoXg%2FIdPc5K9ThBH%2BcgN94Xt8N3f%2BIVTYdzorSig1ONBT2GLTJbtbw84PvW1lzMJ1EDChZEoUXlrDo70l%2FvIG83IdRKzCx0MDkGTDAv4d8RCsbOAEDU9HZspWl%2B9ePS%2BqJsMUCjAR6YmnsoU7AQAa0BpKJOP8M2K2jIJ0VIyw5K2siV7bkhkj%2BhunwNDvx3pfcOTLVtPxL6wEh8E%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2005%3A00%3A04%20GMT&X-Amz-Signature=18a291361ad1216e47bf81fe14ce204faf54f5aa6405dfd236ca2ba1f3824fc6'

$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1750183612496/TAdN8NokD3-camera-video-segment-1750183614996.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSBYDP2TS4%2F20250618%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250618T210004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKr%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIFpWzNZHTIuWx9Jfzil6EW78cqt8sJ%2FyIKXs7AqYf58IAiEA%2FplU9EnLA%2FBKbqOQgPr%2BpbyuqGuuuXTTshUMwZx9cVMq5gMIk%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDGPA%2FyuA9QEJpPr9DCq6Azw6oOyCKDKWrzL%2BHT%2FvK6UaJ2Lzu5Cf4R24oyE%2FlinXD0XsbTuWcmyYsn33d0ziGODpr1IxKv2fGo8pxPMoRS88PhZQOKtx3bXwg9nww1G8NH5neYXf16JZbDbTVVoU%2FAO5jWTR%2BR0%2FscWsGnNq%2FDPDnG%2FDp4ZD4o47mwcRExrnTCRFKs3Bjoi9ZRlko16yUzNitA8Pc8TPl3XgutBSpzQeAnADoQC44bB512exlSR8nRx6vJwvniQv3EAizZwtO8xT8oW%2BU3GkpXLOEa1Ql7u0HXNTl7WRSUm0ZqsLkGHKwxkwkPNy9fCheupmAAewRmeUMF6J4KnidkaH8WdLTp9urm1cFN4Kuutl2npiaOS%2FnhxV%2F3A4SN7XjRZKX89F3uytxQnl%2BjFW9%2BiK%2Ba6tBb34%2FM95IEXlATDB1rHKO5rnuhzqk784KOb7KXNUs76pLZasYyBFfU6s7wVXlO8VQdVI8I56j9JDuOFwqD9Vc%2BpRCp61IkkIp%2BQVIrS0l9VyiFge3JRBXfhU%2BwW94c2VOOSZkCGKotw%2Bp8uzJ9uaXeGYipFxMoLH3JHNgF1BKE42l2TWGb8gjLZBFuswtffLwgY6pQFdUfnEVvfK7XJ%2FoXg%2FIdPc5K9ThBH%2BcgN94Xt8N3f%2BIVTYdzorSig1ONBT2GLTJbtbw84PvW1lzMJ1EDChZEoUXlrDo70l%2FvIG83IdRKzCx0MDkGTDAv4d8RCsbOAEDU9HZspWl%2B9ePS%2BqJsMUCjAR6YmnsoU7AQAa0BpKJOP8M2K2jIJ0VIyw5K2siV7bkhkj%2BhunwNDvx3pfcOTLVtPxL6wEh8E%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2005%3A00%3A04%20GMT&X-Amz-Signature=84827e4ffc3286e76c61b812d29192e612ff189032d5229ddada45ede346e57f'

$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991152736/1750182961683/iBjUYEdqdu-camera-video-segment-driver-1750182964183.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSBYDP2TS4%2F20250618%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250618T210004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKr%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIFpWzNZHTIuWx9Jfzil6EW78cqt8sJ%2FyIKXs7AqYf58IAiEA%2FplU9EnLA%2FBKbqOQgPr%2BpbyuqGuuuXTTshUMwZx9cVMq5gMIk%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDGPA%2FyuA9QEJpPr9DCq6Azw6oOyCKDKWrzL%2BHT%2FvK6UaJ2Lzu5Cf4R24oyE%2FlinXD0XsbTuWcmyYsn33d0ziGODpr1IxKv2fGo8pxPMoRS88PhZQOKtx3bXwg9nww1G8NH5neYXf16JZbDbTVVoU%2FAO5jWTR%2BR0%2FscWsGnNq%2FDPDnG%2FDp4ZD4o47mwcRExrnTCRFKs3Bjoi9ZRlko16yUzNitA8Pc8TPl3XgutBSpzQeAnADoQC44bB512exlSR8nRx6vJwvniQv3EAizZwtO8xT8oW%2BU3GkpXLOEa1Ql7u0HXNTl7WRSUm0ZqsLkGHKwxkwkPNy9fCheupmAAewRmeUMF6J4KnidkaH8WdLTp9urm1cFN4Kuutl2npiaOS%2FnhxV%2F3A4SN7XjRZKX89F3uytxQnl%2BjFW9%2BiK%2Ba6tBb34%2FM95IEXlATDB1rHKO5rnuhzqk784KOb7KXNUs76pLZasYyBFfU6s7wVXlO8VQdVI8I56j9JDuOFwqD9Vc%2BpRCp61IkkIp%2BQVIrS0l9VyiFge3JRBXfhU%2BwW94c2VOOSZkCGKotw%2Bp8uzJ9uaXeGYipFxMoLH3JHNgF1BKE42l2TWGb8gjLZBFuswtffLwgY6pQFdUfnEVvfK7XJ%2FoXg%2FIdPc5K9ThBH%2BcgN94Xt8N3f%2BIVTYdzorSig1ONBT2GLTJbtbw84PvW1lzMJ1EDChZEoUXlrDo70l%2FvIG83IdRKzCx0MDkGTDAv4d8RCsbOAEDU9HZspWl%2B9ePS%2BqJsMUCjAR6YmnsoU7AQAa0BpKJOP8M2K2jIJ0VIyw5K2siV7bkhkj%2BhunwNDvx3pfcOTLVtPxL6wEh8E%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2005%3A00%3A04%20GMT&X-Amz-Signature=a6693641f651b4421729a5556d135c44e208f66c83e36df1319404ecbd0a5de8'

$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991152736/1750182961683/Foph7ZSROm-camera-video-segment-1750182964183.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSBYDP2TS4%2F20250618%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250618T210004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKr%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIFpWzNZHTIuWx9Jfzil6EW78cqt8sJ%2FyIKXs7AqYf58IAiEA%2FplU9EnLA%2FBKbqOQgPr%2BpbyuqGuuuXTTshUMwZx9cVMq5gMIk%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDGPA%2FyuA9QEJpPr9DCq6Azw6oOyCKDKWrzL%2BHT%2FvK6UaJ2Lzu5Cf4R24oyE%2FlinXD0XsbTuWcmyYsn33d0ziGODpr1IxKv2fGo8pxPMoRS88PhZQOKtx3bXwg9nww1G8NH5neYXf16JZbDbTVVoU%2FAO5jWTR%2BR0%2FscWsGnNq%2FDPDnG%2FDp4ZD4o47mwcRExrnTCRFKs3Bjoi9ZRlko16yUzNitA8Pc8TPl3XgutBSpzQeAnADoQC44bB512exlSR8nRx6vJwvniQv3EAizZwtO8xT8oW%2BU3GkpXLOEa1Ql7u0HXNTl7WRSUm0ZqsLkGHKwxkwkPNy9fCheupmAAewRmeUMF6J4KnidkaH8WdLTp9urm1cFN4Kuutl2npiaOS%2FnhxV%2F3A4SN7XjRZKX89F3uytxQnl%2BjFW9%2BiK%2Ba6tBb34%2FM95IEXlATDB1rHKO5rnuhzqk784KOb7KXNUs76pLZasYyBFfU6s7wVXlO8VQdVI8I56j9JDuOFwqD9Vc%2BpRCp61IkkIp%2BQVIrS0l9VyiFge3JRBXfhU%2BwW94c2VOOSZkCGKotw%2Bp8uzJ9uaXeGYipFxMoLH3JHNgF1BKE42l2TWGb8gjLZBFuswtffLwgY6pQFdUfnEVvfK7XJ%2FoXg%2FIdPc5K9ThBH%2BcgN94Xt8N3f%2BIVTYdzorSig1ONBT2GLTJbtbw84PvW1lzMJ1EDChZEoUXlrDo70l%2FvIG83IdRKzCx0MDkGTDAv4d8RCsbOAEDU9HZspWl%2B9ePS%2BqJsMUCjAR6YmnsoU7AQAa0BpKJOP8M2K2jIJ0VIyw5K2siV7bkhkj%2BhunwNDvx3pfcOTLVtPxL6wEh8E%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2005%3A00%3A04%20GMT&X-Amz-Signature=f23742409dae31579b89bf627f84ac0f2af42baa7637f5585a9d8684779c5a34'
